# Updates cryptos list with latest price/volume data (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '65.674.30'
$ws.Cells.Item(2, 5).Value = '  +0.17%  '

$ws.Cells.Item(3, 4).Value = '2.650.31'
$ws.Cells.Item(3, 5).Value = '  -0.30%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '596.90'
$ws.Cells.Item(5, 5).Value = '  -0.05%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '156.64'
$ws.Cells.Item(6, 5).Value = '  +0.49%  '

$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.627'
$ws.Cells.Item(8, 5).Value = '  +4.16%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.127'
$ws.Cells.Item(9, 5).Value = '  +5.00%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.397'
$ws.Cells.Item(10, 5).Value = '  +0.69%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.80'
$ws.Cells.Item(11, 5).Value = '  -1.09%  '

$ws.Cells.Item(12, 5).Value = '  +1.06%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '28.78'
$ws.Cells.Item(13, 5).Value = '  -1.23%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.0000197'
$ws.Cells.Item(14, 5).Value = '  +2.18%  '

$ws.Cells.Item(15, 4).Value = '3.127.38'
$ws.Cells.Item(15, 5).Value = '  -0.17%  '

$ws.Cells.Item(16, 4).Value = '65.495.17'
$ws.Cells.Item(16, 5).Value = '  +0.20%  '

$ws.Cells.Item(17, 4).Value = '2.663.35'
$ws.Cells.Item(17, 5).Value = '  +0.07%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.63'
$ws.Cells.Item(18, 5).Value = '  +1.67%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.74'
$ws.Cells.Item(19, 5).Value = '  -0.71%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.48'
$ws.Cells.Item(20, 5).Value = '  +0.46%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '349.44'
$ws.Cells.Item(21, 5).Value = '  +0.30%  '

$ws.Cells.Item(22, 5).Value = '  +0.26%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '69.02'
$ws.Cells.Item(23, 5).Value = '  -1.16%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.0000113'
$ws.Cells.Item(24, 5).Value = '  +5.91%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.69'
$ws.Cells.Item(25, 5).Value = '  +0.61%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.65'
$ws.Cells.Item(26, 5).Value = '  +1.11%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.59'
$ws.Cells.Item(27, 5).Value = '  +0.26%  '

$ws.Cells.Item(28, 5).Value = '  -1.40%  '

$ws.Cells.Item(29, 5).Value = '  +0.10%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.90'
$ws.Cells.Item(30, 5).Value = '  -1.29%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '530.11'
$ws.Cells.Item(31, 5).Value = '  -1.12%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.11'
$ws.Cells.Item(32, 5).Value = '  -0.57%  '

$ws.Cells.Item(33, 5).Value = '  +1.18%  '

$ws.Cells.Item(34, 2).Value = 'RenderToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '6.43'
$ws.Cells.Item(34, 5).Value = '  -0.56%  '

$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.45'
$ws.Cells.Item(35, 5).Value = '  +1.47%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.421'
$ws.Cells.Item(36, 5).Value = '  +0.32%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '20.34'
$ws.Cells.Item(37, 5).Value = '  +0.34%  '

$ws.Cells.Item(38, 5).Value = '  -0.06%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.93'
$ws.Cells.Item(39, 5).Value = '  -0.01%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '153.14'
$ws.Cells.Item(40, 5).Value = '  -3.34%  '

$ws.Cells.Item(41, 5).Value = '  +0.04%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '160.46'
$ws.Cells.Item(42, 5).Value = '  -2.72%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.07'
$ws.Cells.Item(43, 5).Value = '  +0.65%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.30'
$ws.Cells.Item(44, 5).Value = '  +2.55%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0605'
$ws.Cells.Item(45, 5).Value = '  -0.19%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '22.59'
$ws.Cells.Item(46, 5).Value = '  -0.89%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.634'
$ws.Cells.Item(47, 5).Value = '  -1.66%  '

$ws.Cells.Item(48, 5).Value = '  -0.99%  '

$ws.Cells.Item(49, 4).Value = '0.0₆0255'
$ws.Cells.Item(49, 5).Value = '  +12.86%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0994'
$ws.Cells.Item(50, 5).Value = '  +0.20%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '19.68'
$ws.Cells.Item(51, 5).Value = '  -0.85%  '
